$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 11-14 down to rows 12-15 (First/Last name, Email, Phone, Year columns),
# and insert the new "User8" record at row 11.

$ws.Range("B11").Value = "User8"
$ws.Range("C11").Value = "User8LN"
$ws.Range("D11").Value = "user8@gmail.com"
$ws.Range("E11").Value = ""
$ws.Range("F11").Value = 2000

$ws.Range("B12").Value = "Bruce"
$ws.Range("C12").Value = "Lee"
$ws.Range("D12").Value = "bruce@gmail.com"
$ws.Range("E12").Value = "(111)333445"
$ws.Range("F12").Value = 1987

$ws.Range("B13").Value = "Gamora"
$ws.Range("C13").Value = "Gamorak"
$ws.Range("D13").Value = "gamora@gmail.com"
$ws.Range("E13").Value = "(111)333111"
$ws.Range("F13").Value = 1988

$ws.Range("B14").Value = "Witcher"
$ws.Range("C14").Value = "Moon"
$ws.Range("D14").Value = "witcher@gmail.com"
$ws.Range("E14").Value = "(111)333999"
$ws.Range("F14").Value = 1990

$ws.Range("B15").Value = "Supwom"
$ws.Range("C15").Value = "Nanual"
$ws.Range("D15").Value = "supwom@gmail.com"
$ws.Range("E15").Value = "(111)333777"
$ws.Range("F15").Value = 1988
